# Change the team name textbox on slide 1 from "Team MuffinTime" to "Team Hotdog"
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("TextBox 8")   # id=9, holding the team name

# Replace the whole text (merges the two runs "Team " + "MuffinTime" into one run)
$tr = $shp.TextFrame.TextRange
$full = $tr.Characters(1, $tr.Length)
$full.Text = "Team Hotdog"

# The textbox uses wrap="none" + auto-fit, so PowerPoint re-centers/resizes it
# around the same center point once the (shorter) text is applied.
$shp.Left = 399.5418897637795
$shp.Width = 159.9015
